$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder header row (row 1) ---
# New column order: Date | posWordPercentage | negWordPercentage | posPhrasePercentage |
# negPhrasePercentage | ElapsedMs | wordCount | sentenceCount | posWordCount | negWordCount |
# positivePhraseCount | negativePhraseCount | Method
$ws.Range("B1").Value = "posWordPercentage"
$ws.Range("C1").Value = "negWordPercentage"
$ws.Range("D1").Value = "posPhrasePercentage"
$ws.Range("E1").Value = "negPhrasePercentage"
$ws.Range("F1").Value = "ElapsedMs"
$ws.Range("G1").Value = "wordCount"
$ws.Range("H1").Value = "sentenceCount"
$ws.Range("I1").Value = "posWordCount"
$ws.Range("J1").Value = "negWordCount"
$ws.Range("K1").Value = "positivePhraseCount"
$ws.Range("L1").Value = "negativePhraseCount"
$ws.Range("M1").Value = "Method"

# --- Update data row (row 2) to match new columns / new sample values ---
$ws.Range("A2").Value = 42605.455081018517
$ws.Range("B2").Value = 66
$ws.Range("C2").Value = 31
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = 91
$ws.Range("F2").Value = 4443
$ws.Range("G2").Value = 6623
$ws.Range("H2").Value = 762
$ws.Range("I2").Value = 152
$ws.Range("J2").Value = 71
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 21
$ws.Range("M2").Value = "Bag"

# --- Resize columns to match the new bestFit widths for the reordered headers ---
# (inputs chosen so the engine's ColumnWidth quantization lands as close as
# possible to the exact target OOXML widths from the diff)
$ws.Columns.Item(2).ColumnWidth = 18.5
$ws.Columns.Item(3).ColumnWidth = 18.666666666666668
$ws.Columns.Item(4).ColumnWidth = 19.666666666666668
$ws.Columns.Item(5).ColumnWidth = 19.666666666666668
$ws.Columns.Item(6).ColumnWidth = 9.666666666666666
$ws.Columns.Item(7).ColumnWidth = 10.0
$ws.Columns.Item(8).ColumnWidth = 13.666666666666666
$ws.Columns.Item(9).ColumnWidth = 13.666666666666666
$ws.Columns.Item(10).ColumnWidth = 13.666666666666666
$ws.Columns.Item(11).ColumnWidth = 18.833333333333332
$ws.Columns.Item(12).ColumnWidth = 19.5
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
